$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark checklist rows 5 and 9 (B5, B9) as accomplished with an "X"
$ws.Range("B5").Value = "X"
$ws.Range("B9").Value = "X"

# Move the selection to reflect the author's final cursor position
$ws.Range("B9").Select()
